# EVE bare install tutorial.docx - apply commit "Add files via upload"
#
# Summary of changes applied:
#  1. Replace the "and install Ubuntu server..." paragraph + the UNetLab
#     hyperlink paragraph with two new paragraphs of plain text, the
#     second one carrying the (new) "_GoBack" bookmark at its end.
#  2. Mark "passwd" as a grammar-flagged word (gramStart/gramEnd) and
#     add a new empty "ListParagraph" paragraph right after it.
#  3. Split "Change hostname if needed" into two runs ("...if need" + "ed").
#  4. Mark "nano" (in "nano /etc/hostname and ...") as a grammar-flagged
#     word (gramStart/gramEnd).
#  5. Remove the old "_GoBack" bookmark further down the document (Word
#     relocates it to the edited text automatically on save).
#
# Bookmark w:id numbers are re-assigned by the engine in document order
# at save time, so we don't need to hand compute the +1 shifts the diff
# shows for OLE_LINK4..OLE_LINK9 - they happen automatically once the new
# "_GoBack" bookmark is introduced near the top of the document.

$d = $word.ActiveDocument

# Common OpenXML package wrapper used for every InsertXML call below.
# Bundling a tiny styles.xml part (with just the styles we reference)
# keeps w:rStyle references ("tgc") alive - without it Word silently
# drops rStyle because it can't resolve the style in the mini package.
function New-PkgXml([string]$bodyXml) {
    return @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body></w:document></pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml"><pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="character" w:customStyle="1" w:styleId="tgc"><w:name w:val="_tgc"/><w:basedOn w:val="DefaultParagraphFont"/></w:style><w:style w:type="paragraph" w:styleId="ListParagraph"><w:name w:val="List Paragraph"/></w:style></w:styles></pkg:xmlData></pkg:part>
</pkg:package>
"@
}

function Find-ParagraphByText([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($text + "`r")) {
            return $p
        }
    }
    return $null
}

# --- 0. Remove the stale "_GoBack" bookmark further down the document
#        *before* we introduce the new one below - "_GoBack" is a
#        name that can exist only once, but while both momentarily
#        exist, Bookmarks.Item("_GoBack") resolves to whichever one
#        comes first in the document, so delete the old one first. ---

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 1. Intro paragraphs: replace "and install Ubuntu server..." and the
#        UNetLab hyperlink paragraph with the new two-paragraph intro. ---

$pAnd = Find-ParagraphByText("and install Ubuntu server till you reach UNL installation using bellow guide:")
$pLink = $pAnd.Next()
$introRange = $d.Range($pAnd.Range.Start, $pLink.Range.End)

$introBody = '<w:p><w:r><w:t>Follow the Ubuntu server installation and customize things per your need or leave default.</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">After the installation of Ubuntu </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is done</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, follow the bellow steps in order to install EVE on top.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$introRange.InsertXML((New-PkgXml $introBody))

# --- 2. "passwd: password updated successfully" gets gramStart/gramEnd
#        around "passwd", plus a new empty ListParagraph after it. ---

$pPasswd = Find-ParagraphByText("passwd: password updated successfully")
$passwdBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>passwd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>: password updated successfully</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'

$pPasswd.Range.InsertXML((New-PkgXml $passwdBody))

# --- 3. "Change hostname if need" -> "Change hostname if needed". The
#        diff adds "ed" as its own run rather than merging it into the
#        existing run, so build the paragraph explicitly via InsertXML. ---

$pHostname = Find-ParagraphByText("Change hostname if need")
$hostnameBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Change hostname if need</w:t></w:r><w:r><w:t>ed</w:t></w:r></w:p>'
$pHostname.Range.InsertXML((New-PkgXml $hostnameBody))

# --- 4. "nano /etc/hostname and " gets gramStart/gramEnd around "nano" ---

$pNano = Find-ParagraphByText("nano /etc/hostname and ")
$nanoBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rStyle w:val="tgc"/></w:rPr></w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="OLE_LINK5"/><w:bookmarkStart w:id="1" w:name="OLE_LINK22"/><w:bookmarkStart w:id="2" w:name="OLE_LINK23"/>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t>nano</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t xml:space="preserve"> /</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t>/</w:t></w:r>' + `
    '<w:r><w:rPr><w:rStyle w:val="tgc"/><w:b/><w:bCs/></w:rPr><w:t>hostname</w:t></w:r>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkEnd w:id="1"/><w:bookmarkEnd w:id="2"/>' + `
    '<w:r><w:rPr><w:rStyle w:val="tgc"/></w:rPr><w:t xml:space="preserve">and </w:t></w:r></w:p>'

$pNano.Range.InsertXML((New-PkgXml $nanoBody))

Write-Host "Edit complete."
